$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 1.019
$ws.Range("C2").Value = 1.086
$ws.Range("D2").Value = 1.336
$ws.Range("E2").Value = 0.769
$ws.Range("F2").Value = 1.44

# Row 3 (MSE)
$ws.Range("B3").Value = 1.633
$ws.Range("C3").Value = 1.818
$ws.Range("D3").Value = 3.181
$ws.Range("E3").Value = 1.026
$ws.Range("F3").Value = 4.476

# Row 5 (mean Y-predicted)
$ws.Range("B5").Value = 18.495
$ws.Range("C5").Value = 15.354
$ws.Range("D5").Value = 17.824
$ws.Range("E5").Value = 12.822
$ws.Range("F5").Value = 30.801

# Row 6 (R2)
$ws.Range("B6").Value = 0.769
$ws.Range("C6").Value = 0.868
$ws.Range("D6").Value = 0.842
$ws.Range("E6").Value = 0.45
$ws.Range("F6").Value = 0.886
